$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 1,4
$arr[0,0] = 0.0044
$arr[0,1] = 0.0074
$arr[0,2] = 0.0067
$arr[0,3] = 0.006
$ws.Range("M2:P2").Value = $arr

$ws.Range("R4").Value = -0.028

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0.1538
$arr[0,1] = -0.2559
$arr[0,2] = -0.4623
$arr[0,3] = -0.0669
$arr[0,4] = -0.0832
$arr[0,5] = -0.4215
$arr[0,6] = -0.2976
$arr[0,7] = -0.1307
$arr[0,8] = -0.1664
$arr[0,9] = -1.9431
$ws.Range("M5:V5").Value = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0.2796
$arr[0,1] = 0.16
$arr[0,2] = -0.1199
$arr[0,3] = -0.2378
$arr[0,4] = -0.2232
$arr[0,5] = -0.2835
$arr[0,6] = -0.2191
$arr[0,7] = -0.1703
$arr[0,8] = -0.035
$arr[0,9] = -0.7412
$ws.Range("M6:V6").Value = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = -0.1622
$arr[0,1] = -0.2326
$arr[0,2] = -0.2193
$arr[0,3] = -0.0199
$arr[0,4] = -0.0569
$arr[0,5] = 0.0123
$arr[0,6] = 0.0299
$arr[0,7] = 0.0425
$arr[0,8] = 0.056
$arr[0,9] = 0.0702
$ws.Range("M7:V7").Value = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0.3659
$arr[0,1] = 0.2453
$arr[0,2] = 0.2073
$arr[0,3] = 0.3473
$arr[0,4] = 0.2554
$arr[0,5] = -0.0867
$arr[0,6] = -0.1053
$arr[0,7] = -0.079
$arr[0,8] = -0.073
$arr[0,9] = -0.0034
$ws.Range("M8:V8").Value = $arr

$ws.Range("P10").Value = -0.0433

$ws.Range("R10").Value = -0.0103

$ws.Range("M11").Value = 0.0098

$arr = New-Object "object[,]" 1,10
$arr[0,0] = -0.1481
$arr[0,1] = -0.173
$arr[0,2] = -0.0803
$arr[0,3] = -0.036
$arr[0,4] = -0.0298
$arr[0,5] = -0.0563
$arr[0,6] = -0.0513
$arr[0,7] = -0.0506
$arr[0,8] = -0.105
$arr[0,9] = -6.3047
$ws.Range("M12:V12").Value = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = -0.1408
$arr[0,1] = -0.1035
$arr[0,2] = -0.2906
$arr[0,3] = -0.2615
$arr[0,4] = -0.1048
$arr[0,5] = -0.057
$arr[0,6] = -0.0917
$arr[0,7] = -0.0677
$arr[0,8] = -0.0298
$arr[0,9] = -0.0778
$ws.Range("M14:V14").Value = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0.2312
$arr[0,1] = -0.2678
$arr[0,2] = -0.6669
$arr[0,3] = -0.4807
$arr[0,4] = -0.2571
$arr[0,5] = -0.8043
$arr[0,6] = -0.8209
$arr[0,7] = -0.5283
$arr[0,8] = -0.3389
$arr[0,9] = -50.7598
$ws.Range("M16:V16").Value = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 27957
$arr[0,1] = 28276.7244
$arr[0,2] = 28614.6055
$arr[0,3] = 28959.3736
$arr[0,4] = 29316.3505
$arr[0,5] = 29680.1101
$arr[0,6] = 30035.2088
$arr[0,7] = 30385.2987
$arr[0,8] = 30722.8667
$arr[0,9] = 31060.6435
$ws.Range("M17:V17").Value = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0.1455
$arr[0,1] = -0.119
$arr[0,2] = -0.3516
$arr[0,3] = -0.0106
$arr[0,4] = 0.0153
$arr[0,5] = -0.3047
$arr[0,6] = -0.1895
$arr[0,7] = -0.0363
$arr[0,8] = -0.0874
$arr[0,9] = -1.6966
$ws.Range("M18:V18").Value = $arr

$arr = New-Object "object[,]" 1,4
$arr[0,0] = -0.1019
$arr[0,1] = -0.0998
$arr[0,2] = -0.0967
$arr[0,3] = -0.0937
$ws.Range("P20:S20").Value = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0.2484
$arr[0,1] = 0.0996
$arr[0,2] = 0.3612
$arr[0,3] = -0.0437
$arr[0,4] = -0.0625
$arr[0,5] = 0.2327
$arr[0,6] = 0.1212
$arr[0,7] = -0.0292
$arr[0,8] = 0.0161
$arr[0,9] = -41.6959
$ws.Range("M22:V22").Value = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = -0.027
$arr[0,1] = -0.016
$arr[0,2] = -0.0121
$arr[0,3] = -0.0046
$arr[0,4] = 0.0121
$arr[0,5] = 0.0263
$arr[0,6] = 0.0059
$arr[0,7] = 0.006
$arr[0,8] = 0.011
$arr[0,9] = 0.0157
$ws.Range("M23:V23").Value = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0.0574
$arr[0,1] = 0.1506
$arr[0,2] = 0.0366
$arr[0,3] = -0.0311
$arr[0,4] = 0.0947
$arr[0,5] = -0.0475
$arr[0,6] = -0.0287
$arr[0,7] = -0.0114
$arr[0,8] = -0.0048
$arr[0,9] = 0.0019
$ws.Range("M24:V24").Value = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0.1388
$arr[0,1] = 0.1693
$arr[0,2] = 0.163
$arr[0,3] = 0.1135
$arr[0,4] = 0.0774
$arr[0,5] = 0.0358
$arr[0,6] = -0.0382
$arr[0,7] = -0.0102
$arr[0,8] = -0.0014
$arr[0,9] = -0.0023
$ws.Range("M25:V25").Value = $arr

$ws.Range("M26").Value = 0.0016

$arr = New-Object "object[,]" 1,7
$arr[0,0] = -0.0013
$arr[0,1] = 0.0043
$arr[0,2] = 0.0061
$arr[0,3] = 0.0057
$arr[0,4] = 0.0059
$arr[0,5] = 0.0061
$arr[0,6] = 0.0049
$ws.Range("N29:T29").Value = $arr

$arr = New-Object "object[,]" 1,4
$arr[0,0] = -0.0001
$arr[0,1] = 0.0028
$arr[0,2] = 0.0015
$arr[0,3] = 0.0005
$ws.Range("M30:P30").Value = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0.4485
$arr[0,1] = -0.0006
$arr[0,2] = 0.0003
$arr[0,3] = -0.0003
$arr[0,4] = -0.0003
$arr[0,5] = -0.0003
$arr[0,6] = -0.0003
$arr[0,7] = -0.0004
$arr[0,8] = -0.0004
$arr[0,9] = 0.0003
$ws.Range("M33:V33").Value = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = -0.0008
$arr[0,1] = -0.0038
$arr[0,2] = -0.0017
$arr[0,3] = 0.0002
$arr[0,4] = 0.0011
$arr[0,5] = 0.0011
$arr[0,6] = 0.001
$arr[0,7] = 0.001
$arr[0,8] = 0.0009
$arr[0,9] = 0.0012
$ws.Range("M34:V34").Value = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0.0003
$arr[0,1] = 0.0296
$arr[0,2] = 0.0378
$arr[0,3] = 0.0497
$arr[0,4] = 0.0503
$arr[0,5] = 0.0233
$arr[0,6] = 0.0163
$arr[0,7] = 0.0056
$arr[0,8] = 0.0045
$arr[0,9] = 0.0034
$ws.Range("M35:V35").Value = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0.0004
$arr[0,1] = 0.0155
$arr[0,2] = 0.0262
$arr[0,3] = 0.022
$arr[0,4] = 0.0163
$arr[0,5] = 0.0146
$arr[0,6] = 0.0144
$arr[0,7] = 0.0141
$arr[0,8] = 0.0137
$arr[0,9] = 0.0068
$ws.Range("M36:V36").Value = $arr

$ws.Range("M39").Value = -0.0003

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0.0063
$arr[0,1] = -0.0115
$arr[0,2] = -0.0209
$arr[0,3] = -0.0239
$arr[0,4] = -0.0296
$arr[0,5] = -0.0123
$arr[0,6] = -0.0032
$arr[0,7] = -0.0002
$arr[0,8] = -0.0002
$arr[0,9] = -0.0027
$ws.Range("M40:V40").Value = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = -0.0001
$arr[0,1] = -0.0001
$arr[0,2] = -0.0001
$arr[0,3] = -0.0002
$arr[0,4] = -0.0002
$arr[0,5] = -0.0002
$arr[0,6] = -0.0002
$arr[0,7] = -0.0002
$arr[0,8] = -0.0002
$arr[0,9] = -0.0001
$ws.Range("M42:V42").Value = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0.0443
$arr[0,1] = 0.1365
$arr[0,2] = 0.1324
$arr[0,3] = 0.1589
$arr[0,4] = 0.1171
$arr[0,5] = 0.0157
$arr[0,6] = 0.0154
$arr[0,7] = 0.0039
$arr[0,8] = 0.0295
$arr[0,9] = -0.1477
$ws.Range("M44:V44").Value = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 12.4
$arr[0,1] = 12.5418
$arr[0,2] = 12.6917
$arr[0,3] = 12.8446
$arr[0,4] = 13.0029
$arr[0,5] = 13.1643
$arr[0,6] = 13.3218
$arr[0,7] = 13.477
$arr[0,8] = 13.6268
$arr[0,9] = 13.7766
$ws.Range("M45:V45").Value = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0.4433
$arr[0,1] = -0.0007
$arr[0,2] = 0.0002
$arr[0,3] = -0.0003
$arr[0,4] = -0.0003
$arr[0,5] = -0.0003
$arr[0,6] = -0.0004
$arr[0,7] = -0.0005
$arr[0,8] = -0.0005
$arr[0,9] = 0.0001
$ws.Range("M46:V46").Value = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = -0.3824
$arr[0,1] = 0.0535
$arr[0,2] = -0.0002
$arr[0,3] = 0.0213
$arr[0,4] = 0.0318
$arr[0,5] = 0.0001
$arr[0,6] = 0.0001
$arr[0,7] = 0.0002
$arr[0,8] = 0.0002
$arr[0,9] = -0.1618
$ws.Range("M50:V50").Value = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = -0.0091
$arr[0,1] = 0.0003
$arr[0,2] = 0.0004
$arr[0,3] = 0.0004
$arr[0,4] = 0.0004
$arr[0,5] = 0.0004
$arr[0,6] = 0.0004
$arr[0,7] = 0.0004
$arr[0,8] = 0.0004
$arr[0,9] = 0.0004
$ws.Range("M51:V51").Value = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0.0405
$arr[0,1] = 0.1015
$arr[0,2] = 0.1043
$arr[0,3] = 0.1076
$arr[0,4] = 0.068
$arr[0,5] = 0.0089
$arr[0,6] = 0.0055
$arr[0,7] = 0.0015
$arr[0,8] = 0.0016
$arr[0,9] = 0.0017
$ws.Range("M52:V52").Value = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = -0.0594
$arr[0,1] = -0.0462
$arr[0,2] = -0.0109
$arr[0,3] = -0.0151
$arr[0,4] = -0.0177
$arr[0,5] = -0.0185
$arr[0,6] = -0.0183
$arr[0,7] = -0.0181
$arr[0,8] = 0.0091
$arr[0,9] = 0.0032
$ws.Range("M53:V53").Value = $arr

$ws.Range("M54").Value = 0.0004

$arr = New-Object "object[,]" 1,8
$arr[0,0] = -0.0002
$arr[0,1] = -0.0002
$arr[0,2] = -0.0003
$arr[0,3] = -0.0002
$arr[0,4] = -0.0002
$arr[0,5] = -0.0001
$arr[0,6] = -0.0001
$arr[0,7] = -0.0001
$ws.Range("N57:U57").Value = $arr

